$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (e.g. "0.602", "1.00").
# Excel auto-converts Range.Value assignments that parse as numbers into numeric
# cells, which would lose the original text formatting (trailing zeros, etc.) that
# the source data stores as literal text. Force text interpretation per-cell first,
# then restore the default (unstyled) look after the value is written, since the
# diff shows no style/numFmt change on any of these cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply updated cell values (Coin / Link / Price / Volume(1h))
$ws.Range("D2").Value = '67.054.09'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '3.533.64'
$ws.Range("E3").Value = '  +1.44%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '589.35'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").Value = '177.85'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.529.87'
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").Value = '6.96'
$ws.Range("E11").Value = '  -1.57%  '
$ws.Range("D12").Value = '0.426'
$ws.Range("E12").Value = '  -2.11%  '
$ws.Range("D13").Value = '4.140.70'
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").Value = '30.65'
$ws.Range("E14").Value = '  -3.98%  '
$ws.Range("E15").Value = '  -2.32%  '
$ws.Range("D16").Value = '67.001.76'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").Value = '3.539.78'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '6.13'
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").Value = '14.10'
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").Value = '384.86'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("D22").Value = '7.90'
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = '0.543'
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '5.76'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").Value = '71.99'
$ws.Range("E26").Value = '  -2.51%  '
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  -3.39%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '24.56'
$ws.Range("E31").Value = '  +4.48%  '
$ws.Range("D32").Value = '5.97'
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("E34").Value = '  -3.47%  '
$ws.Range("D35").Value = '7.29'
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").Value = '1.59'
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = '29.52'
$ws.Range("E38").Value = '  +13.00%  '
$ws.Range("D39").Value = '159.95'
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("D40").Value = '0.896'
$ws.Range("E40").Value = '  +3.15%  '
$ws.Range("D41").Value = '1.82'
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("E42").Value = '  -2.56%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.61'
$ws.Range("E43").Value = '  -5.21%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '4.55'
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.744.83'
$ws.Range("E45").Value = '  -3.12%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '0.0712'
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("D47").Value = '25.65'
$ws.Range("E47").Value = '  -5.10%  '
$ws.Range("D48").Value = '40.79'
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").Value = '327.53'
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("E51").Value = '  -1.80%  '

# Restore default (no explicit number format) style on the cells forced to text above
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
